$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44 / Row 46 swap: ImmutableX <-> MantraDAO ---
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"

# --- D column (price) updates: force text format for cells that would otherwise parse as numbers ---
$textForceCells = @("D5", "D6", "D7", "D8", "D10", "D12", "D13", "D14", "D18", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D43", "D48", "D49", "D50", "D51", "D44", "D46")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "93.329.10"
$ws.Range("D3").Value = "3.416.48"
$ws.Range("D5").Value = "233.16"
$ws.Range("D6").Value = "620.45"
$ws.Range("D7").Value = "1.43"
$ws.Range("D8").Value = "0.392"
$ws.Range("D10").Value = "0.974"
$ws.Range("D11").Value = "3.415.16"
$ws.Range("D12").Value = "43.17"
$ws.Range("D13").Value = "0.199"
$ws.Range("D14").Value = "6.29"
$ws.Range("D15").Value = "93.191.70"
$ws.Range("D16").Value = "4.052.42"
$ws.Range("D18").Value = "8.23"
$ws.Range("D19").Value = "3.414.82"
$ws.Range("D20").Value = "18.07"
$ws.Range("D22").Value = "0.494"
$ws.Range("D23").Value = "3.39"
$ws.Range("D24").Value = "497.89"
$ws.Range("D25").Value = "6.65"
$ws.Range("D26").Value = "0.0000184"
$ws.Range("D27").Value = "95.54"
$ws.Range("D28").Value = "12.01"
$ws.Range("D29").Value = "3.599.77"
$ws.Range("D30").Value = "11.36"
$ws.Range("D32").Value = "0.138"
$ws.Range("D34").Value = "0.995"
$ws.Range("D35").Value = "0.174"
$ws.Range("D36").Value = "0.548"
$ws.Range("D37").Value = "28.98"
$ws.Range("D38").Value = "566.22"
$ws.Range("D39").Value = "7.48"
$ws.Range("D40").Value = "1.41"
$ws.Range("D43").Value = "0.898"
$ws.Range("D48").Value = "5.46"
$ws.Range("D49").Value = "53.35"
$ws.Range("D50").Value = "2.12"
$ws.Range("D51").Value = "8.14"
$ws.Range("D44").Value = "3.72"
$ws.Range("D46").Value = "1.70"

foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}

# --- E column (volume %) updates ---
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +2.53%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  +6.64%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("E20").Value = "  +6.26%  "
$ws.Range("E21").Value = "  +5.98%  "
$ws.Range("E22").Value = "  +8.74%  "
$ws.Range("E23").Value = "  +7.41%  "
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("E27").Value = "  +6.34%  "
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +4.82%  "
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("E36").Value = "  +2.82%  "
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("E38").Value = "  +5.84%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("E47").Value = "  +4.61%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("E51").Value = "  +3.58%  "
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("E46").Value = "  +2.47%  "
